$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.643.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.98%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.826.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.33%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'309.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.64%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.22%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4664"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.24%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07149"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.00%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9049"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.38%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07701"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.53%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'19.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.826.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.12%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.272"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.375"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.73%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'87.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.29%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008566"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.13%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'26.680.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.36%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.032"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.15%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.908"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.20%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.25%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.78%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.988"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.83%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'113.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.78%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.872"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.08816"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.72%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.832"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +5.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7351"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.441"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.11%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.35%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.01934"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.24%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.881"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.64%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5070"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.27%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1498"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.56%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.078"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.84%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.008"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.22%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4671"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.95%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.20%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'98.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.57%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.576"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.06037"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.27%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'64.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'35.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.87%  "
$ws.Range("E51").Style = "Normal"
